$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# "Picture 2" (id=1026) - shape index 9 on this slide
$pic2 = $s.Shapes.Item(9)
$pic2.Left = 483.2815748031496
$pic2.Top = 300.27134708267715
$pic2.Width = 205.3552017503937
$pic2.Height = 51.85212598425197

# "Picture 5" (id=6) - shape index 11 on this slide
$pic5 = $s.Shapes.Item(11)
$pic5.Left = 729.25
$pic5.Top = 70.43354420708661
$pic5.Width = 195.90251968503938
$pic5.Height = 78.8507882015748
